$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column E header, copying the header style/format from D1 first
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Plaza Vea Universitaria"

# Update existing data rows 2-5 with the new multi-format values and fill column E
$ws.Range("B2").Value = 9702.98
$ws.Range("C2").Value = 3798.9
$ws.Range("D2").Value = 706.93
$ws.Range("E2").Value = 442.12

$ws.Range("B3").Value = 16816.150000000001
$ws.Range("C3").Value = 8029.09
$ws.Range("D3").Value = 3896.83
$ws.Range("E3").Value = 1840.42

$ws.Range("B4").Value = 615.92999999999995
$ws.Range("C4").Value = 202.91
$ws.Range("D4").Value = 61.33
$ws.Range("E4").Value = 41.54

$ws.Range("B5").Value = 3744.17
$ws.Range("C5").Value = 3758.41
$ws.Range("D5").Value = 3346.19
$ws.Range("E5").Value = 879.84

# Rows 6-15 keep zero values in B,C,D but now also need a zero in column E
for ($r = 6; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = 0
}

# Update the visible selection to include the new column E
$ws.Range("A1:E15").Select()

$wb.Save()
